$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 8
$ws.Cells.Item(8, 8).Value = 89.72727
$ws.Cells.Item(8, 9).Value = 89.72727
$ws.Cells.Item(8, 11).Value = 269.18181
$ws.Cells.Item(8, 13).Value = -130.18181

# ALC row 53
$ws.Cells.Item(53, 8).Value = 279.83334
$ws.Cells.Item(53, 9).Value = 163.58333
$ws.Cells.Item(53, 11).Value = 163.58333
$ws.Cells.Item(53, 13).Value = 473.41667

# ALC row 86
$ws.Cells.Item(86, 8).Value = 2581.7273
$ws.Cells.Item(86, 9).Value = 2585.25
$ws.Cells.Item(86, 11).Value = 2585.25
$ws.Cells.Item(86, 13).Value = -1462.25

# ALC row 89
$ws.Cells.Item(89, 8).Value = 2581.7273
$ws.Cells.Item(89, 9).Value = 2585.25
$ws.Cells.Item(89, 11).Value = 12926.25
$ws.Cells.Item(89, 13).Value = -7310.25

# ALC row 103
$ws.Cells.Item(103, 8).Value = 1862.091
$ws.Cells.Item(103, 9).Value = 779.8
$ws.Cells.Item(103, 11).Value = 2339.4
$ws.Cells.Item(103, 13).Value = -1753.4

# ALC row 106
$ws.Cells.Item(106, 8).Value = 1499.5
$ws.Cells.Item(106, 9).Value = 1499.5
$ws.Cells.Item(106, 11).Value = 1499.5
$ws.Cells.Item(106, 13).Value = -868.5

# ALC row 112
$ws.Cells.Item(112, 8).Value = 2553.158
$ws.Cells.Item(112, 10).Value = 2653.647
$ws.Cells.Item(112, 12).Value = 7960.941
$ws.Cells.Item(112, 14).Value = -10176.941

# ALC row 132
$ws.Cells.Item(132, 8).Value = 3281.1936
$ws.Cells.Item(132, 9).Value = 2321.05
$ws.Cells.Item(132, 11).Value = 6963.150000000001
$ws.Cells.Item(132, 13).Value = -4433.150000000001

# ALC row 137
$ws.Cells.Item(137, 8).Value = 2209.5217
$ws.Cells.Item(137, 10).Value = 4950
$ws.Cells.Item(137, 12).Value = 14850
$ws.Cells.Item(137, 14).Value = -19950

# ALC row 138
$ws.Cells.Item(138, 8).Value = 3602.647
$ws.Cells.Item(138, 10).Value = 4228.143
$ws.Cells.Item(138, 12).Value = 12684.429
$ws.Cells.Item(138, 14).Value = -22964.429

$ws = $wb.Worksheets.Item("ARM")
# ARM row 63
$ws.Cells.Item(63, 8).Value = 4516.4546
$ws.Cells.Item(63, 9).Value = 2837.2
$ws.Cells.Item(63, 10).Value = 5915.8335
$ws.Cells.Item(63, 11).Value = 2837.2
$ws.Cells.Item(63, 12).Value = 5915.8335
$ws.Cells.Item(63, 13).Value = -2151.2
$ws.Cells.Item(63, 14).Value = -7287.8335

# ARM row 66
$ws.Cells.Item(66, 8).Value = 4516.4546
$ws.Cells.Item(66, 9).Value = 2837.2
$ws.Cells.Item(66, 10).Value = 5915.8335
$ws.Cells.Item(66, 11).Value = 14186
$ws.Cells.Item(66, 12).Value = 29579.1675
$ws.Cells.Item(66, 13).Value = -10754
$ws.Cells.Item(66, 14).Value = -36443.1675

# ARM row 88
$ws.Cells.Item(88, 8).Value = 4946.4546
$ws.Cells.Item(88, 10).Value = 5450.875
$ws.Cells.Item(88, 12).Value = 5450.875
$ws.Cells.Item(88, 14).Value = -6262.875

# ARM row 91
$ws.Cells.Item(91, 8).Value = 4946.4546
$ws.Cells.Item(91, 10).Value = 5450.875
$ws.Cells.Item(91, 12).Value = 5450.875
$ws.Cells.Item(91, 14).Value = -8258.875

# ARM row 132
$ws.Cells.Item(132, 8).Value = 1320.25
$ws.Cells.Item(132, 9).Value = 1320.25
$ws.Cells.Item(132, 11).Value = 3960.75
$ws.Cells.Item(132, 13).Value = -1430.75

$ws = $wb.Worksheets.Item("BSM")
# BSM row 82
$ws.Cells.Item(82, 8).Value = 5128.5
$ws.Cells.Item(82, 9).Value = 5128.5
$ws.Cells.Item(82, 11).Value = 5128.5
$ws.Cells.Item(82, 13).Value = -4745.5

# BSM row 85
$ws.Cells.Item(85, 8).Value = 5128.5
$ws.Cells.Item(85, 9).Value = 5128.5
$ws.Cells.Item(85, 11).Value = 5128.5
$ws.Cells.Item(85, 13).Value = -3802.5

# BSM row 97
$ws.Cells.Item(97, 8).Value = 17466
$ws.Cells.Item(97, 9).Value = 12945.6
$ws.Cells.Item(97, 10).Value = 25000
$ws.Cells.Item(97, 11).Value = 12945.6
$ws.Cells.Item(97, 12).Value = 25000
$ws.Cells.Item(97, 13).Value = -11954.6
$ws.Cells.Item(97, 14).Value = -26982

$ws = $wb.Worksheets.Item("CRP")
# CRP row 6
$ws.Cells.Item(6, 8).Value = 666742
$ws.Cells.Item(6, 9).Value = 800050.4
$ws.Cells.Item(6, 10).Value = 200
$ws.Cells.Item(6, 11).Value = 800050.4
$ws.Cells.Item(6, 12).Value = 200
$ws.Cells.Item(6, 13).Value = -799937.4
$ws.Cells.Item(6, 14).Value = -426

# CRP row 31
$ws.Cells.Item(31, 8).Value = 2619.111
$ws.Cells.Item(31, 9).Value = 2759.25
$ws.Cells.Item(31, 10).Value = 1498
$ws.Cells.Item(31, 11).Value = 2759.25
$ws.Cells.Item(31, 12).Value = 1498
$ws.Cells.Item(31, 13).Value = -2464.25
$ws.Cells.Item(31, 14).Value = -2088

# CRP row 34
$ws.Cells.Item(34, 8).Value = 2619.111
$ws.Cells.Item(34, 9).Value = 2759.25
$ws.Cells.Item(34, 10).Value = 1498
$ws.Cells.Item(34, 11).Value = 2759.25
$ws.Cells.Item(34, 12).Value = 1498
$ws.Cells.Item(34, 13).Value = -2557.25
$ws.Cells.Item(34, 14).Value = -1902

# CRP row 59
$ws.Cells.Item(59, 8).Value = 49998
$ws.Cells.Item(59, 10).Value = 49998
$ws.Cells.Item(59, 12).Value = 49998
$ws.Cells.Item(59, 14).Value = -52288

# CRP row 105
$ws.Cells.Item(105, 8).Value = 903.75
$ws.Cells.Item(105, 9).Value = 872.3333
$ws.Cells.Item(105, 10).Value = 998
$ws.Cells.Item(105, 11).Value = 872.3333
$ws.Cells.Item(105, 12).Value = 998
$ws.Cells.Item(105, 13).Value = 874.6667
$ws.Cells.Item(105, 14).Value = -4492

# CRP row 114
$ws.Cells.Item(114, 8).Value = 99500
$ws.Cells.Item(114, 10).Value = 99500
$ws.Cells.Item(114, 12).Value = 99500
$ws.Cells.Item(114, 14).Value = -108178

$ws = $wb.Worksheets.Item("CUL")
# CUL row 70
$ws.Cells.Item(70, 8).Value = 2393.5
$ws.Cells.Item(70, 9).Value = 524.6667
$ws.Cells.Item(70, 11).Value = 1574.0001
$ws.Cells.Item(70, 13).Value = -1259.0001

# CUL row 73
$ws.Cells.Item(73, 8).Value = 2393.5
$ws.Cells.Item(73, 9).Value = 524.6667
$ws.Cells.Item(73, 11).Value = 1574.0001
$ws.Cells.Item(73, 13).Value = -482.0001

# CUL row 131
$ws.Cells.Item(131, 8).Value = 1131.125
$ws.Cells.Item(131, 9).Value = 1049.5
$ws.Cells.Item(131, 11).Value = 3148.5
$ws.Cells.Item(131, 13).Value = 1891.5

$ws = $wb.Worksheets.Item("GSM")
# GSM row 49
$ws.Cells.Item(49, 8).Value = 42500
$ws.Cells.Item(49, 10).Value = 42500
$ws.Cells.Item(49, 12).Value = 42500
$ws.Cells.Item(49, 14).Value = -42868

# GSM row 52
$ws.Cells.Item(52, 8).Value = 25000
$ws.Cells.Item(52, 10).Value = 25000
$ws.Cells.Item(52, 12).Value = 25000
$ws.Cells.Item(52, 14).Value = -25518

# GSM row 70
$ws.Cells.Item(70, 8).Value = 12502475
$ws.Cells.Item(70, 9).Value = 14288114
$ws.Cells.Item(70, 11).Value = 14288114
$ws.Cells.Item(70, 13).Value = -14287844

# GSM row 73
$ws.Cells.Item(73, 8).Value = 12502475
$ws.Cells.Item(73, 9).Value = 14288114
$ws.Cells.Item(73, 11).Value = 14288114
$ws.Cells.Item(73, 13).Value = -14287178

$ws = $wb.Worksheets.Item("LTW")
# LTW row 19
$ws.Cells.Item(19, 8).Value = 2249
$ws.Cells.Item(19, 9).Value = 1999
$ws.Cells.Item(19, 10).Value = 2499
$ws.Cells.Item(19, 11).Value = 1999
$ws.Cells.Item(19, 12).Value = 2499
$ws.Cells.Item(19, 13).Value = -1829
$ws.Cells.Item(19, 14).Value = -2839

# LTW row 22
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(22, 8).Value = 1650
$ws.Cells.Item(22, 9).Value = 1650
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 1650
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -1355

# LTW row 27
$ws.Cells.Item(27, 14).ClearContents()
$ws.Cells.Item(27, 8).Value = 1650
$ws.Cells.Item(27, 9).Value = 1650
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 11).Value = 1650
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 13).Value = -1543

# LTW row 132
$ws.Cells.Item(132, 8).Value = 5246.278
$ws.Cells.Item(132, 9).Value = 5433.375
$ws.Cells.Item(132, 10).Value = 3749.5
$ws.Cells.Item(132, 11).Value = 16300.125
$ws.Cells.Item(132, 12).Value = 11248.5
$ws.Cells.Item(132, 13).Value = -13770.125
$ws.Cells.Item(132, 14).Value = -16308.5

$ws = $wb.Worksheets.Item("WVR")
# WVR row 55
$ws.Cells.Item(55, 8).Value = 4499.2
$ws.Cells.Item(55, 10).Value = 5499
$ws.Cells.Item(55, 12).Value = 5499
$ws.Cells.Item(55, 14).Value = -6053

# WVR row 132
$ws.Cells.Item(132, 8).Value = 804.35
$ws.Cells.Item(132, 9).Value = 583.5263
$ws.Cells.Item(132, 11).Value = 1750.5789
$ws.Cells.Item(132, 13).Value = 779.4211
